$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 0.6666666666666666
$ws.Range("G2").Value = 0.08785
$ws.Range("H2").Value = 0.26355
$ws.Range("M2").Value = 7.163958
$ws.Range("N2").Value = 21.491874
$ws.Range("O2").Value = 0.5236853002292368
$ws.Range("P2").Value = 0.5236853002292368
$ws.Range("Q2").Value = 0.6293537103
$ws.Range("R2").Value = 5.6641833927
$ws.Range("S2").Value = 0.5236853002292368
$ws.Range("T2").Value = 0.5236853002292368

# Row 3
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 0.6666666666666666
$ws.Range("G3").Value = 0.08785
$ws.Range("H3").Value = 0.26355
$ws.Range("N3").Value = 6.273242999999999
$ws.Range("O3").Value = 0.1528580124686175
$ws.Range("P3").Value = 0.1528580124686176
$ws.Range("Q3").Value = 0.1837014658499999
$ws.Range("R3").Value = 1.65331319265
$ws.Range("S3").Value = 0.1528580124686175
$ws.Range("T3").Value = 0.1528580124686176

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.08785
$ws.Range("H4").Value = 0.26355
$ws.Range("M4").Value = 0.4156183333333334
$ws.Range("N4").Value = 1.246855
$ws.Range("O4").Value = 0.03038169845111343
$ws.Range("P4").Value = 0.03038169845111344
$ws.Range("Q4").Value = 0.03651207058333333
$ws.Range("R4").Value = 0.32860863525
$ws.Range("S4").Value = 0.03038169845111343
$ws.Range("T4").Value = 0.03038169845111344

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.08785
$ws.Range("H5").Value = 0.26355
$ws.Range("M5").Value = 3.342326666666667
$ws.Range("N5").Value = 10.02698
$ws.Range("O5").Value = 0.2443240655371678
$ws.Range("P5").Value = 0.2443240655371678
$ws.Range("Q5").Value = 0.2936233976666667
$ws.Range("R5").Value = 2.642610579
$ws.Range("S5").Value = 0.2443240655371678
$ws.Range("T5").Value = 0.2443240655371678

# Row 6
$ws.Range("E6").Value = 2
$ws.Range("F6").Value = 0.6666666666666666
$ws.Range("G6").Value = 0.08785
$ws.Range("H6").Value = 0.26355
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 0.1557673333333333
$ws.Range("N6").Value = 0.467302
$ws.Range("O6").Value = 0.01138659142370381
$ws.Range("P6").Value = 0.01138659142370381
$ws.Range("Q6").Value = 0.01368416023333333
$ws.Range("R6").Value = 0.1231574421
$ws.Range("S6").Value = 0.01138659142370381
$ws.Range("T6").Value = 0.01138659142370381

# Row 7
$ws.Range("E7").Value = 2
$ws.Range("F7").Value = 0.6666666666666666
$ws.Range("G7").Value = 0.08785
$ws.Range("H7").Value = 0.26355
$ws.Range("M7").Value = 0.51114
$ws.Range("N7").Value = 1.53342
$ws.Range("O7").Value = 0.03736433189016073
$ws.Range("P7").Value = 0.03736433189016074
$ws.Range("Q7").Value = 0.044903649
$ws.Range("R7").Value = 0.404132841
$ws.Range("S7").Value = 0.03736433189016073
$ws.Range("T7").Value = 0.03736433189016074
